# Updated estimates - slipped about a week (actually two weeks) on the
# three URI-related items, whose drafts moved from AlexJ to Alex.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("schedule")

# Author re-assigned from "AlexJ" to "Alex" for the URI items (rows 19-21).
$ws.Range("B19").Value = "Alex"
$ws.Range("B20").Value = "Alex"
$ws.Range("B21").Value = "Alex"

# Draft ETA slipped two weeks for the same three items.
$ws.Range("C19").Value = 40977
$ws.Range("C20").Value = 40981
$ws.Range("C21").Value = 40984

# Reflect where the author was last working in the sheet.
$ws.Range("D37").Select()
